$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 7.235148
$ws.Range("H2").Value = 21.705444
$ws.Range("I2").Value = 0.9254344869740032
$ws.Range("J2").Value = 0.9254344869740032
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 4.847498666666667
$ws.Range("N2").Value = 14.542496
$ws.Range("O2").Value = 0.03400671694637637
$ws.Range("P2").Value = 0.03400671694637637
$ws.Range("Q2").Value = 35.072370283136
$ws.Range("R2").Value = 315.651332548224
$ws.Range("S2").Value = 0.03147098865093995
$ws.Range("T2").Value = 0.03147098865093995

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 7.235148
$ws.Range("H3").Value = 21.705444
$ws.Range("I3").Value = 0.9254344869740032
$ws.Range("J3").Value = 0.9254344869740032
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.627093333333334
$ws.Range("N3").Value = 4.88128
$ws.Range("O3").Value = 0.01141456784970118
$ws.Range("P3").Value = 0.01141456784970118
$ws.Range("Q3").Value = 11.77226107648
$ws.Range("R3").Value = 105.95034968832
$ws.Range("S3").Value = 0.01056343474201816
$ws.Range("T3").Value = 0.01056343474201816

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 7.235148
$ws.Range("H4").Value = 21.705444
$ws.Range("I4").Value = 0.9254344869740032
$ws.Range("J4").Value = 0.9254344869740032
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 136.0707373333333
$ws.Range("N4").Value = 408.212212
$ws.Range("O4").Value = 0.9545787152039225
$ws.Range("P4").Value = 0.9545787152039225
$ws.Range("Q4").Value = 984.491923075792
$ws.Range("R4").Value = 8860.427307682128
$ws.Range("S4").Value = 0.8834000635810451
$ws.Range("T4").Value = 0.8834000635810451

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.3016356666666667
$ws.Range("H5").Value = 0.9049070000000001
$ws.Range("I5").Value = 0.03858166390441884
$ws.Range("J5").Value = 0.03858166390441884
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 4.847498666666667
$ws.Range("N5").Value = 14.542496
$ws.Range("O5").Value = 0.03400671694637637
$ws.Range("P5").Value = 0.03400671694637637
$ws.Range("Q5").Value = 1.462178491985778
$ws.Range("R5").Value = 13.159606427872
$ws.Range("S5").Value = 0.001312035723717798
$ws.Range("T5").Value = 0.001312035723717798

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.3016356666666667
$ws.Range("H6").Value = 0.9049070000000001
$ws.Range("I6").Value = 0.03858166390441884
$ws.Range("J6").Value = 0.03858166390441884
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.627093333333334
$ws.Range("N6").Value = 4.88128
$ws.Range("O6").Value = 0.01141456784970118
$ws.Range("P6").Value = 0.01141456784970118
$ws.Range("Q6").Value = 0.490789382328889
$ws.Range("R6").Value = 4.417104440960001
$ws.Range("S6").Value = 0.0004403930203913558
$ws.Range("T6").Value = 0.0004403930203913558

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.3016356666666667
$ws.Range("H7").Value = 0.9049070000000001
$ws.Range("I7").Value = 0.03858166390441884
$ws.Range("J7").Value = 0.03858166390441884
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 136.0707373333333
$ws.Range("N7").Value = 408.212212
$ws.Range("O7").Value = 0.9545787152039225
$ws.Range("P7").Value = 0.9545787152039225
$ws.Range("Q7").Value = 41.04378756936489
$ws.Range("R7").Value = 369.3940881242841
$ws.Range("S7").Value = 0.03682923516030968
$ws.Range("T7").Value = 0.03682923516030968

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.2813256666666666
$ws.Range("H8").Value = 0.843977
$ws.Range("I8").Value = 0.0359838491215779
$ws.Range("J8").Value = 0.0359838491215779
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 4.847498666666667
$ws.Range("N8").Value = 14.542496
$ws.Range("O8").Value = 0.03400671694637637
$ws.Range("P8").Value = 0.03400671694637637
$ws.Range("Q8").Value = 1.363725794065778
$ws.Range("R8").Value = 12.273532146592
$ws.Range("S8").Value = 0.001223692571718613
$ws.Range("T8").Value = 0.001223692571718613

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.2813256666666666
$ws.Range("H9").Value = 0.843977
$ws.Range("I9").Value = 0.0359838491215779
$ws.Range("J9").Value = 0.0359838491215779
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.627093333333334
$ws.Range("N9").Value = 4.88128
$ws.Range("O9").Value = 0.01141456784970118
$ws.Range("P9").Value = 0.01141456784970118
$ws.Range("Q9").Value = 0.4577431167288889
$ws.Range("R9").Value = 4.11968805056
$ws.Range("S9").Value = 0.0004107400872916612
$ws.Range("T9").Value = 0.0004107400872916612

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.2813256666666666
$ws.Range("H10").Value = 0.843977
$ws.Range("I10").Value = 0.0359838491215779
$ws.Range("J10").Value = 0.0359838491215779
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 136.0707373333333
$ws.Range("N10").Value = 408.212212
$ws.Range("O10").Value = 0.9545787152039225
$ws.Range("P10").Value = 0.9545787152039225
$ws.Range("Q10").Value = 38.28019089412489
$ws.Range("R10").Value = 344.521718047124
$ws.Range("S10").Value = 0.03434941646256762
$ws.Range("T10").Value = 0.03434941646256762

